$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column G
$ws.Range("G1").Value = "expected status code"

# New "expected status code" values for rows 2-6 (all 422 - invalid credential cases)
$ws.Range("G2:G6").Value = 422

# Match the number format used by column F (style index 2 -> numFmtId 49 "@")
$ws.Range("G1:G6").NumberFormat = "@"

# Last scenario (row 6) now also exercises the "yes" (Email) case like row 2/row 4
$ws.Range("A6").Value = "yes"

# New column G needs the same width treatment as the others (stored width ends up
# 5/6 character wider than the ColumnWidth we set, so compensate to land on 23)
$ws.Columns.Item(7).ColumnWidth = 23 - (5/6)
